# Weekly update: shift the "Fruta / hortaliza" Piña price history down by one
# week (4 rows per week: Especial / Primera / Segunda / Tercera) and insert a
# new, most-recent week at the top of the data block (row 261).
#
# Inserting whole rows above row 261 pushes the existing rows 261:296 down to
# 265:300 (which is exactly the body of the diff - each former row r becomes
# r+4), grows the used range to A1:T300, and leaves a blank 4-row gap at
# 261:264 ready to be populated with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("261:264").Insert()

# Columns that stay constant across the whole "Agrícola del Norte S.A. de
# Arica" / Piña / Caramelo block, for every quality row of the new week.
$common = @{
    A = 1
    B = "Agrícola del Norte S.A. de Arica"
    C = "Arica y Parinacota"
    D = 45013
    E = 15
    F = "Fruta"
    G = 100108
    H = "Tropicales y subtropicales"
    I = 100108005
    J = "Piña"
    K = "Caramelo"
    R = "Ecuador"
}

# Per-quality-row figures for the new week (row offset 0..3 from row 261).
$weekRows = @(
    @{ L = "Especial"; M = 200; N = 20000; O = 22000; P = 21000; Q = "`$/caja 10 unidades"; S = 2100; T = 10 }
    @{ L = "Primera";  M = 250; N = 20000; O = 22000; P = 21000; Q = "`$/caja 12 unidades"; S = 1750; T = 12 }
    @{ L = "Segunda";  M = 260; N = 20000; O = 22000; P = 21000; Q = "`$/caja 14 unidades"; S = 1500; T = 14 }
    @{ L = "Tercera";  M = 200; N = 20000; O = 22000; P = 21000; Q = "`$/caja 16 unidades"; S = 1312; T = 16 }
)

for ($i = 0; $i -lt $weekRows.Count; $i++) {
    $r = 261 + $i
    $row = $weekRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $common.A
    $ws.Cells.Item($r, 2).Value2 = $common.B
    $ws.Cells.Item($r, 3).Value2 = $common.C
    $ws.Cells.Item($r, 4).Value2 = $common.D
    $ws.Cells.Item($r, 5).Value2 = $common.E
    $ws.Cells.Item($r, 6).Value2 = $common.F
    $ws.Cells.Item($r, 7).Value2 = $common.G
    $ws.Cells.Item($r, 8).Value2 = $common.H
    $ws.Cells.Item($r, 9).Value2 = $common.I
    $ws.Cells.Item($r, 10).Value2 = $common.J
    $ws.Cells.Item($r, 11).Value2 = $common.K
    $ws.Cells.Item($r, 12).Value2 = $row.L
    $ws.Cells.Item($r, 13).Value2 = $row.M
    $ws.Cells.Item($r, 14).Value2 = $row.N
    $ws.Cells.Item($r, 15).Value2 = $row.O
    $ws.Cells.Item($r, 16).Value2 = $row.P
    $ws.Cells.Item($r, 17).Value2 = $row.Q
    $ws.Cells.Item($r, 18).Value2 = $common.R
    $ws.Cells.Item($r, 19).Value2 = $row.S
    $ws.Cells.Item($r, 20).Value2 = $row.T
}
